# Update cryptocurrency price/volume data per upstream refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.626.92'
$ws.Range('E2').Value = '  +0.11%  '
$ws.Range('D3').Value = '2.648.48'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''604.23'
$ws.Range('D6').Value = '''146.84'
$ws.Range('E6').Value = '  +1.90%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +0.43%  '
$ws.Range('E9').Value = '  +1.14%  '
$ws.Range('E10').Value = '  -0.37%  '
$ws.Range('E11').Value = '  +4.44%  '
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('D14').Value = '3.125.93'
$ws.Range('E14').Value = '  +0.10%  '
$ws.Range('D15').Value = '63.480.17'
$ws.Range('E15').Value = '  +0.02%  '
$ws.Range('E16').Value = '  +0.86%  '
$ws.Range('D17').Value = '2.655.95'
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('D18').Value = '''11.51'
$ws.Range('E18').Value = '  +0.60%  '
$ws.Range('E19').Value = '  +4.38%  '
$ws.Range('D20').Value = '''343.23'
$ws.Range('E20').Value = '  +0.79%  '
$ws.Range('E21').Value = '  +2.80%  '
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').Value = '''5.59'
$ws.Range('E23').Value = '  -3.09%  '
$ws.Range('E24').Value = '  -0.47%  '
$ws.Range('E25').Value = '  +1.41%  '
$ws.Range('E26').Value = '  +7.64%  '
$ws.Range('D27').Value = '''573.86'
$ws.Range('E27').Value = '  +4.83%  '
$ws.Range('E28').Value = '  +1.03%  '
$ws.Range('E29').Value = '  -1.56%  '
$ws.Range('D30').Value = '''7.99'
$ws.Range('E30').Value = '  +2.84%  '
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('D32').Value = '''2.04'
$ws.Range('E32').Value = '  +3.63%  '
$ws.Range('D34').Value = '0.0₃0822'
$ws.Range('E34').Value = '  +1.73%  '
$ws.Range('E35').Value = '  +6.61%  '
$ws.Range('D36').Value = '''168.64'
$ws.Range('E36').Value = '  -3.78%  '
$ws.Range('E37').Value = '  +1.30%  '
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range('E39').Value = '  +7.11%  '
$ws.Range('D40').Value = '''19.11'
$ws.Range('E40').Value = '  +0.13%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').Value = '''169.13'
$ws.Range('E42').Value = '  -0.63%  '
$ws.Range('D44').Value = '''22.24'
$ws.Range('E44').Value = '  -0.65%  '
$ws.Range('E45').Value = '  +2.74%  '
$ws.Range('E46').Value = '  +0.17%  '
$ws.Range('E47').Value = '  +3.10%  '
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').Value = '''1.89'
$ws.Range('E49').Value = '  +10.72%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''18.85'
$ws.Range('E50').Value = '  +0.17%  '
$ws.Range('E51').Value = '  +1.94%  '
